$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 29) following the same pattern as the existing
# "historico*" rows (26-28): column A = lowercase entity name,
# column B = "<Entity>Controller", column C = "SP_<ENTITY>_SELECT".
# Values are written C, B, A to match the shared-string insertion order
# of the original edit.
$ws.Range("C29").Value = "SP_BONOSHIS_RESUMEN_SELECT"
$ws.Range("B29").Value = "BonoshisresumenController"
$ws.Range("A29").Value = "bonoshisresumen"

# Copy the formatting used by the previous row so the new row matches
# the existing style (fill/border) of the table.
$ws.Range("A28:C28").Copy()
$ws.Range("A29:C29").PasteSpecial(-4122)

# Update the active selection to the newly added cell, matching the
# author's last interaction with the sheet.
$ws.Range("A29").Select()
